# Apply the stadium quiz edits:
#  - Column A (type) gets a value for rows 2-15 (=1) and rows 16-25 (=3)
#  - Column F (slugCorAns) gets overwritten to "no_pic" for rows 2-25

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 15; $r++) {
    $ws.Cells.Item($r, 1).Value = 1
}

for ($r = 16; $r -le 25; $r++) {
    $ws.Cells.Item($r, 1).Value = 3
}

for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 6).Value = "no_pic"
}
